# Applies the reshuffled species-observation data for rows 2-7.
# Net effect: the data rows (everything that varies per-observation) get
# cyclically permuted: new row 2 <- old row 3, new row 3 <- old row 7,
# new row 4 <- old row 5, new row 5 <- old row 6, new row 6 <- old row 4,
# new row 7 <- old row 2. Columns that are identical across all the rows
# (C, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY) are left
# untouched since their values don't change either way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param($Row, $A, $B, $D, $E, $F, $G, $H, $Q, $R)
    $ws.Cells.Item($Row, 1).Value = $A    # A - Id
    $ws.Cells.Item($Row, 2).Value = $B    # B - Taxonsorteringsordning
    $ws.Cells.Item($Row, 4).Value = $D    # D - Rödlistade
    $ws.Cells.Item($Row, 5).Value = $E    # E - TaxonId
    $ws.Cells.Item($Row, 6).Value = $F    # F - Artnamn
    $ws.Cells.Item($Row, 7).Value = $G    # G - Vetenskapligt namn
    $ws.Cells.Item($Row, 8).Value = $H    # H - Auktor
    $ws.Cells.Item($Row, 17).Value = $Q   # Q - Ost
    $ws.Cells.Item($Row, 18).Value = $R   # R - Nord
}

# Row 2 <- old row 3 data
Set-RowData 2 111402342 77267 "NT" 6446 "Kolflarnlav" "Carbonicola anthracophila" "(Nyl.) Bendiksby & Timdal" 545197.7500373307 7020179.372318991

# Row 3 <- old row 7 data (also gains K/L/M/N = activity data)
Set-RowData 3 111402340 56414 "NT" 100049 "Spillkråka" "Dryocopus martius" "(Linnaeus, 1758)" 545198.1129081531 7020057.514641967
$ws.Cells.Item(3, 11).Value = ""          # K3
$ws.Cells.Item(3, 12).Value = ""          # L3
$ws.Cells.Item(3, 13).Value = "spel/sång" # M3 - Aktivitet
$ws.Cells.Item(3, 14).Value = ""          # N3

# Row 4 <- old row 5 data
Set-RowData 4 111402343 77186 "NT" 353 "Dvärgbägarlav" "Cladonia parasitica" "(Hoffm.) Hoffm." 545197.7500373307 7020179.372318991

# Row 5 <- old row 6 data
Set-RowData 5 111402339 77267 "NT" 6446 "Kolflarnlav" "Carbonicola anthracophila" "(Nyl.) Bendiksby & Timdal" 545198.1129081531 7020057.514641967

# Row 6 <- old row 4 data
Set-RowData 6 111402344 90666 "LC" 4364 "Dropptaggsvamp" "Hydnellum ferrugineum" "(Fr.:Fr.) P. Karst." 545197.7500373307 7020179.372318991

# Row 7 <- old row 2 data (loses K/L/M/N activity data)
Set-RowData 7 111402337 96370 "LC" 219847 "Tvåblad" "Neottia ovata" "(L.) Buff. & Fingerh." 545198.1129081531 7020057.514641967
$ws.Cells.Item(7, 11).ClearContents()  # K7
$ws.Cells.Item(7, 12).ClearContents()  # L7
$ws.Cells.Item(7, 13).ClearContents()  # M7
$ws.Cells.Item(7, 14).ClearContents()  # N7
